$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.885879039764404
$ws.Range("B1").Value = 4.103427410125732
$ws.Range("C1").Value = 3.388366222381592
$ws.Range("D1").Value = 2.870473146438599
$ws.Range("E1").Value = 1.295594692230225
